# Finish project - update student IDs in the "For students" checklist rows
# 24-26 and move the active selection/tab from "For graders" to
# "For students" (cell A27:D27 selected, scrolled to A3).

$wb = $excel.ActiveWorkbook

$wsStudents = $wb.Worksheets.Item("For students")

# --- Update the mis-filled Student ID cells (rows 24-26) on "For students" ---
$wsStudents.Range("C24").Value2 = 23127216
$wsStudents.Range("D24").Value2 = 23127223

$wsStudents.Range("C25").Value2 = 23127216
$wsStudents.Range("D25").Value2 = 23127223
$wsStudents.Range("F25").Value2 = 23127223

$wsStudents.Range("C26").Value2 = 23127216
$wsStudents.Range("D26").Value2 = 23217223

# --- Switch the active sheet/selection to "For students" ---
$wsStudents.Activate()
$wsStudents.Range("A27:D27").Select()

$excel.ActiveWindow.ScrollRow = 3
